# Test_input_Condition2Pos.xlsx - "tiny change to excel input for cond 2"
#
# - Rename sheets: Sheet3 -> lipids, Sheet4 -> files
# - Lowercase the three header labels that moved sheets/got re-typed upstream:
#     "Name" -> "name", "Formula" -> "formula" (lipids sheet, row 1)
#     "Filename" -> "filename" (files sheet, row 1)
# - Leave every other header/value untouched.
# - Selection/active-sheet ends up on the "files" sheet with A2 selected,
#   and the "lipids" sheet ends with H1 selected (no longer the tab shown).

$wb = $excel.ActiveWorkbook

$lipids = $wb.Worksheets.Item(1)
$files  = $wb.Worksheets.Item(2)

$lipids.Name = "lipids"
$files.Name  = "files"

$lipids.Range("A1").Value = "name"
$lipids.Range("B1").Value = "formula"
$files.Range("A1").Value  = "filename"

# Final on-screen selection state: lipids!H1 selected (not the active tab),
# files!A2 selected and files is the active/front sheet.
$null = $lipids.Range("H1").Select()
$null = $files.Activate()
$null = $files.Range("A2").Select()
